# LOB1215.docx restructuring: sections get reshuffled so that the
# "Docente(s) Responsavel(eis)" / "Programa resumido" / "Programa" blocks
# move ahead of the "Objetivos" paragraphs, the evaluation paragraph gains
# the previously-separate "O metodo..." sentence, and the bibliography
# content slides up into the "Avaliacao" bullet while "Bibliografia" +
# the docente line move to the very end of the document.
#
# Strategy: every paragraph keeps its original style/position in the
# Paragraphs collection (verified against the target XML) - only the
# *text* inside each paragraph needs to change. So we read the current
# text out of each paragraph first (to avoid retyping long accented
# strings) and then write it back into its new home with Find/Execute,
# scoped to a single paragraph's Range so there is no cross-paragraph
# ambiguity.

$d = $word.ActiveDocument

function Get-ParaText($idx) {
    $t = $d.Paragraphs.Item($idx).Range.Text
    return $t.Substring(0, $t.Length - 1)
}

function Replace-InPara($idx, $old, $new) {
    $r = $d.Paragraphs.Item($idx).Range
    $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# ---- capture the original content we will need to relocate ----------
$txtFornecer = Get-ParaText 6    # "Fornecer aos alunos..." (PT objective)
$txtProvide  = Get-ParaText 7    # "Provide to students..." (EN objective, italic)
$txtAnaKarine = Get-ParaText 9   # "7043088 - Ana Karine Furtado de Carvalho"
$txtEnergyShort = Get-ParaText 12 # "Energy sources...regional development" (EN, italic)
$txtBiblio = Get-ParaText 19     # full bibliography block

# ---- paragraph 6: becomes the Portuguese "Recursos energéticos..." ---
# (currently duplicated at paragraph 11, which keeps its own copy)
$txtRecursos1 = Get-ParaText 11
Replace-InPara 6 $txtFornecer $txtRecursos1

# ---- paragraph 7: becomes the short English paragraph (was at 12) ---
Replace-InPara 7 $txtProvide $txtEnergyShort

# ---- paragraph 9 (ListBullet): becomes the "Fornecer aos alunos..." --
Replace-InPara 9 $txtAnaKarine $txtFornecer

# ---- paragraph 11: gains the extra closing sentence ------------------
Replace-InPara 11 "no desenvolvimento regional." "no desenvolvimento regiona. A disciplina pode contar com viagens didáticas para complementação do conteúdo da disciplina."

# ---- paragraph 12: becomes the original "Provide to students..." ----
Replace-InPara 12 $txtEnergyShort $txtProvide

# ---- paragraph 14: becomes the short "O método de avaliação..." line -
Replace-InPara 14 "no desenvolvimento regiona. A disciplina pode contar com viagens didáticas para complementação do conteúdo da disciplina." "O método de avaliação será composto por avaliação teórica, apresentação escrita e oral."

# ---- paragraph 17: rotate the three answer segments, then append ----
# the bibliography text after "Norma de recuperação: ". Order matters:
# go from the last segment to the first so no intermediate text becomes
# ambiguous (duplicated) before it is consumed.
Replace-InPara 17 "Avaliação de recuperação (R) envolvendo todo o conteúdo da disciplina. Média Final = (NF+R) / 2 => 5,0 Aprovado" $txtBiblio
Replace-InPara 17 "Para o cálculo da nota final (NF) será adotada a média ponderada de provas e atividades." "Avaliação de recuperação (R) envolvendo todo o conteúdo da disciplina. Média Final = (NF+R) / 2 => 5,0 Aprovado"
Replace-InPara 17 "O método de avaliação será composto por avaliação teórica, apresentação escrita e oral." "Para o cálculo da nota final (NF) será adotada a média ponderada de provas e atividades."

# ---- paragraph 19: becomes the docente line that used to be at 9 ----
Replace-InPara 19 $txtBiblio $txtAnaKarine

Write-Output "done"
